# B6-PowerPoint.pptx edit
#
# 1) Three tables (on the slides that hold them) get switched from the
#    deck's custom "Table_0" style to the built-in "Medium Style 2 -
#    Accent 1" table style.
# 2) The colour scheme that is actually painted on every slide (carried
#    by the sole slide master / theme1.xml, currently the "Integral"
#    theme's "Red Violet" palette) is switched to the stock "Office"
#    palette.

function HexToRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------
$newTableStyleId = "{3E9D2CFE-A63B-4371-A47E-FDA32039C848}"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Swap the slide master's colour scheme to the "Office" palette
$officeHex = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$colorScheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = HexToRGB $officeHex[$i - 1]
}
